$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-26"

# Update the header label for the April column
$ws.Range("A5").Value = "April (through 04-26)"

# Update April (row 5) values
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 54
$ws.Range("E5").Value = 46
$ws.Range("F5").Value = 39
$ws.Range("G5").Value = 53
$ws.Range("H5").Value = 90
$ws.Range("I5").Value = 108

# Update Total (row 6) values
$ws.Range("B6").Value = 83
$ws.Range("C6").Value = 158
$ws.Range("D6").Value = 243
$ws.Range("E6").Value = 243
$ws.Range("F6").Value = 149
$ws.Range("G6").Value = 251
$ws.Range("H6").Value = 513
$ws.Range("I6").Value = 543
